# "excel hz sprint 1 update"
# Mark two Sprint1 stories (US06 / US10) as Done, fill in their actual
# size/time, completed flag, and source/test traceability columns, add two
# retro notes, and leave the Sprint1 sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# ---- Row 7 : US06 / divorceBeforeDeath -----------------------------------
$ws.Range("D7").Value = "Done"
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = 15
$ws.Range("I7").Value = $true

$ws.Range("K7").Value = "hzSprint1.py"
$ws.Range("K7").HorizontalAlignment = -4131

$ws.Range("L7").Value = "divorceBeforeDeath"

$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = "5-12"
$ws.Range("M7").HorizontalAlignment = -4131

$ws.Range("O7").Value = "hzSprint1.py"
$ws.Range("O7").HorizontalAlignment = -4131

$ws.Range("P7").Value = "test_US06.py"
$ws.Range("P7").HorizontalAlignment = -4131

$ws.Range("Q7").Value = "39-49"
$ws.Range("Q7").HorizontalAlignment = -4131

# ---- Row 9 : US10 / marriageAfterAge --------------------------------------
$ws.Range("D9").Value = "Done"
$ws.Range("G9").Value = 12
$ws.Range("H9").Value = 20
$ws.Range("I9").Value = $true

$ws.Range("K9").Value = "hzSprint1.py"
$ws.Range("K9").HorizontalAlignment = -4131

$ws.Range("L9").Value = "marriageAfterAge"

$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "15-27"
$ws.Range("M9").HorizontalAlignment = -4131

$ws.Range("O9").Value = "hzSprint1.py"
$ws.Range("O9").HorizontalAlignment = -4131

$ws.Range("P9").Value = "test_US10.py"
$ws.Range("P9").HorizontalAlignment = -4131

$ws.Range("Q9").Value = "51-61"
$ws.Range("Q9").HorizontalAlignment = -4131

# ---- Retro notes -----------------------------------------------------------
$ws.Range("B14").Value = "Object-Oriented design"
$ws.Range("B14").WrapText = $true

$ws.Range("B18").Value = "Messy github directory"
$ws.Range("B18").WrapText = $true

# ---- Make Sprint1 the active sheet / selection -----------------------------
$ws.Activate()
$ws.Range("D21").Select()
